# Generate Report for Handoff
# Updates the "Status" column (shared across sheets) from "In Translation"
# to "Ready for handoff", and refreshes the handoff timestamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status column: "In Translation" -> "Ready for handoff"
# (shared string used by Overview!B2, Overview!C2, zh-cn!C2, de-de!C2)
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Latest Handoff Date on Overview sheet
$overview.Range("D2").Value = "2016-26-18 10:26:35"

# Latest Handoff Datetime on the per-language sheets
$zhcn.Range("E2").Value = "2016-03-18 10:26:32"
$dede.Range("E2").Value = "2016-03-18 10:26:35"
